$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($range, [string]$val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-CellText $ws.Range("D2") '27.446.83'
Set-CellText $ws.Range("E2") '  +4.64%  '
Set-CellText $ws.Range("D3") '1.719.79'
Set-CellText $ws.Range("E3") '  +3.84%  '
Set-CellText $ws.Range("D4") '1.000'
Set-CellText $ws.Range("E4") '  -0.34%  '
Set-CellText $ws.Range("D5") '230.44'
Set-CellText $ws.Range("E5") '  +4.90%  '
Set-CellText $ws.Range("D6") '0.5437'
Set-CellText $ws.Range("E6") '  +4.02%  '
Set-CellText $ws.Range("E7") '  -0.40%  '
Set-CellText $ws.Range("D8") '0.2774'
Set-CellText $ws.Range("E8") '  +4.31%  '
Set-CellText $ws.Range("D9") '0.06513'
Set-CellText $ws.Range("E9") '  +2.86%  '
Set-CellText $ws.Range("D10") '21.70'
Set-CellText $ws.Range("E10") '  +5.16%  '
Set-CellText $ws.Range("D11") '0.07759'
Set-CellText $ws.Range("E11") '  -0.10%  '
Set-CellText $ws.Range("D12") '4.739'
Set-CellText $ws.Range("E12") '  +4.19%  '
Set-CellText $ws.Range("D13") '1.733.62'
Set-CellText $ws.Range("E13") '  +4.68%  '
Set-CellText $ws.Range("D14") '1.954.35'
Set-CellText $ws.Range("E14") '  +3.67%  '
Set-CellText $ws.Range("D15") '0.6039'
Set-CellText $ws.Range("E15") '  +6.36%  '
Set-CellText $ws.Range("D16") '0.0₅8332'
Set-CellText $ws.Range("E16") '  +2.70%  '
Set-CellText $ws.Range("D17") '69.18'
Set-CellText $ws.Range("E17") '  +5.67%  '
Set-CellText $ws.Range("D18") '27.384.28'
Set-CellText $ws.Range("E18") '  +4.38%  '
Set-CellText $ws.Range("D19") '4.833'
Set-CellText $ws.Range("E19") '  +2.23%  '
Set-CellText $ws.Range("D20") '211.32'
Set-CellText $ws.Range("E20") '  +9.69%  '
Set-CellText $ws.Range("D21") '1.001'
Set-CellText $ws.Range("E21") '  -0.35%  '
Set-CellText $ws.Range("D22") '11.01'
Set-CellText $ws.Range("E22") '  +6.61%  '
Set-CellText $ws.Range("D23") '6.264'
Set-CellText $ws.Range("E23") '  +3.84%  '
Set-CellText $ws.Range("D24") '1.001'
Set-CellText $ws.Range("E24") '  -0.37%  '
Set-CellText $ws.Range("D25") '147.22'
Set-CellText $ws.Range("E25") '  +2.42%  '
Set-CellText $ws.Range("D26") '0.1258'
Set-CellText $ws.Range("E26") '  +4.84%  '
Set-CellText $ws.Range("D27") '7.432'
Set-CellText $ws.Range("E27") '  +2.14%  '
Set-CellText $ws.Range("D28") '17.11'
Set-CellText $ws.Range("E28") '  +6.91%  '
Set-CellText $ws.Range("D29") '1.626'
Set-CellText $ws.Range("E29") '  +8.65%  '
Set-CellText $ws.Range("D30") '0.05664'
Set-CellText $ws.Range("E30") '  +0.99%  '
Set-CellText $ws.Range("D31") '1.316'
Set-CellText $ws.Range("E31") '  +2.73%  '
Set-CellText $ws.Range("D32") '3.658'
Set-CellText $ws.Range("E32") '  +4.40%  '
Set-CellText $ws.Range("D33") '3.518'
Set-CellText $ws.Range("E33") '  +4.10%  '
Set-CellText $ws.Range("E34") '  +3.45%  '
Set-CellText $ws.Range("D35") '0.9775'
Set-CellText $ws.Range("E35") '  +3.45%  '
Set-CellText $ws.Range("D37") '2.430'
Set-CellText $ws.Range("E37") '  +1.00%  '
Set-CellText $ws.Range("D38") '0.5852'
Set-CellText $ws.Range("E38") '  +1.44%  '
Set-CellText $ws.Range("D39") '0.01641'
Set-CellText $ws.Range("E39") '  +2.67%  '
Set-CellText $ws.Range("D40") '5.899'
Set-CellText $ws.Range("E40") '  -0.24%  '
Set-CellText $ws.Range("D41") '1.050.87'
Set-CellText $ws.Range("E41") '  +1.73%  '
Set-CellText $ws.Range("D42") '0.9994'
Set-CellText $ws.Range("E42") '  -0.46%  '
Set-CellText $ws.Range("D43") '0.8368'
Set-CellText $ws.Range("E43") '  -1.35%  '
Set-CellText $ws.Range("D44") '102.20'
Set-CellText $ws.Range("E44") '  -0.11%  '
Set-CellText $ws.Range("D45") '1.859.85'
Set-CellText $ws.Range("E45") '  +3.59%  '
Set-CellText $ws.Range("D46") '60.07'
Set-CellText $ws.Range("E46") '  +2.81%  '
Set-CellText $ws.Range("E47") '  +6.02%  '
Set-CellText $ws.Range("D48") '8.132'
Set-CellText $ws.Range("E48") '  +1.72%  '
Set-CellText $ws.Range("D49") '0.9932'
Set-CellText $ws.Range("E49") '  -0.86%  '
Set-CellText $ws.Range("D50") '0.4349'
Set-CellText $ws.Range("E50") '  -0.10%  '
Set-CellText $ws.Range("D51") '0.05263'
Set-CellText $ws.Range("E51") '  -1.00%  '
